$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-23 17:19:47"
$wsZhCn.Range("H2").Value = "2016-03-23 17:20:13"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-23 17:19:52"
$wsDeDe.Range("H2").Value = "2016-03-23 17:20:19"
